# Apply the "elbasvir resistance data" update:
#  - Insert a new data row for elbasvir (NS5A inhibitor, Merck, MK-8742)
#  - Re-sort the whole table by category (NS3/4A, NS5A x3, NS3/4A, NS5B)
#  - Normalise the non-breaking space in the category text to a plain space
#  - Carry the "bottom of table" marker border down to the new last row
#  - Leave a couple of incidental formatted-but-empty cells behind (B8, A16:B16)
#  - Nudge the saved window position

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a row so the new elbasvir entry has somewhere to live --------
$ws.Rows.Item(3).Insert()

# --- Final table contents (already resorted by category) -----------------
$data = @(
    @("id",           "abbreviation", "category",                        "producer",         "researchCode", "featureRequiringCoverage"),
    @("glecaprevir",   "GLE",          "NS3/4A protease inhibitors",      "Abbvie",           "ABT-493",      "NS3"),
    @("elbasvir",      "EBR",          "NS5A inhibitors",                 "Merck",            "MK-8742",      "NS5A"),
    @("velpatasvir",   "VEL",          "NS5A inhibitors",                 "Gilead Sciences",  "GS-5816",      "NS5A"),
    @("pibrentasvir",  "PIB",          "NS5A inhibitors",                 "Abbvie",           "ABT-530",      "NS5A"),
    @("voxilaprevir",  "VOX",          "NS3/4A protease inhibitors",      "Gilead Sciences",  "GS-9857",      "NS3"),
    @("sofosbuvir",    "SOF",          "NS5B RNA polymerase inhibitors",  "Gilead Sciences",  "GS-7977",      "NS5B")
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $row = $data[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($r + 1, $c + 1).Value = $row[$c]
    }
}

# --- The old "last row" marker style (border flag, no visible border) now
#     belongs to the new last data row (sofosbuvir), not voxilaprevir ------
$ws.Range("B5").ClearFormats()
$ws.Range("B7").Style = "Normal"
$ws.Range("B7").Borders.Item(9).LineStyle = -4142

# --- Incidental formatting left on a couple of empty cells below the table
$blank = $ws.Range("B8")
$blank.Font.Name = "Calibri"
$blank.Font.Size = 11
$blank.Borders.Item(7).LineStyle = 1
$blank.Borders.Item(8).LineStyle = 1
$blank.Borders.Item(9).LineStyle = 1
$blank.HorizontalAlignment = -4131

$ws.Range("A16:B16").Borders.Item(9).LineStyle = -4142
$ws.Range("A16:B16").Interior.Pattern = -4142

# --- Selection / view bookkeeping -----------------------------------------
$ws.Range("A11").Select()
$excel.ActiveWindow.Left = 21620
$excel.ActiveWindow.Top = 9960
